# The roster table's "Name" column (column 4) needs every entry shifted
# down by one row: row 2 (Burugula, Karthik) is cleared out, and every
# other name moves into the next row down, with the previously-empty
# row 10 picking up "Shriram" (what used to be in row 9).
#
# Word's spell-check leaves <w:proofErr> bookmarks glued to specific
# words, so the exact run/proofErr layout that lands in each cell is
# reproduced below (taken from the corresponding source cell) rather
# than simply re-typing the name text.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$pPr = '<w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0563C1"/><w:u w:val="single"/></w:rPr></w:pPr>'

function Set-NameCell($row, $paraId, $textId, $inner) {
    $cell = $t.Cell($row, 4)
    $p = '<w:p w14:paraId="' + $paraId + '" w14:textId="' + $textId + '" w:rsidR="006C710E" w:rsidRDefault="006C710E">' + $pPr + $inner + '</w:p>'
    $xml = $pkgHeader + $p + $pkgFooter
    $cell.Range.InsertXML($xml)
}

$runProps = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0563C1"/><w:u w:val="single"/></w:rPr>'

function Make-Run($text) {
    return '<w:r>' + $runProps + '<w:t>' + $text + '</w:t></w:r>'
}

function Make-SpellCheckedRun($text) {
    return '<w:proofErr w:type="spellStart"/>' + (Make-Run $text) + '<w:proofErr w:type="spellEnd"/>'
}

# Row 2 (paraId 5F6608AD): "Burugula, Karthik (KABURUGU)" -> now empty.
Set-NameCell 2 "5F6608AD" "334D53DE" ""

# Row 3 (paraId 6DB39B73): "Chinmay, Nandi (NCHINMAY)" -> "Burugula, Karthik (KABURUGU)"
Set-NameCell 3 "6DB39B73" "2456275C" ((Make-SpellCheckedRun "Burugula") + (Make-Run ", Karthik (KABURUGU)"))

# Row 4 (paraId 51C2982C): "Praneetha" -> "Chinmay, Nandi (NCHINMAY)"
Set-NameCell 4 "51C2982C" "300B4FDD" (Make-Run "Chinmay, Nandi (NCHINMAY)")

# Row 5 (paraId 1914E90A): "Linga" -> "Praneetha"
Set-NameCell 5 "1914E90A" "38734473" (Make-SpellCheckedRun "Praneetha")

# Row 6 (paraId 07BC27FD): "Abhiram" -> "Linga"
Set-NameCell 6 "07BC27FD" "4136F74F" (Make-Run "Linga")

# Row 7 (paraId 13D0F576): "Thrisundar" -> "Abhiram"
Set-NameCell 7 "13D0F576" "77656C82" (Make-SpellCheckedRun "Abhiram")

# Row 8 (paraId 159B646A): "Uttam" -> "Thrisundar"
Set-NameCell 8 "159B646A" "32312F8B" (Make-SpellCheckedRun "Thrisundar")

# Row 9 (paraId 7D2C69DD): "Shriram" -> "Uttam"
Set-NameCell 9 "7D2C69DD" "4799582A" (Make-Run "Uttam")

# Row 10 (paraId 3BA3C4FD): empty -> "Shriram"
Set-NameCell 10 "3BA3C4FD" "77777777" (Make-Run "Shriram")
